# Active_Outages.xlsx update — 6/19/2025, 9:15:20 AM
# - Refresh the "Elapsed Duration(Hrs)" text for every currently-open outage
#   (time has moved on by 1 minute 44 seconds since the previous snapshot).
# - Append a newly-detected outage row (JED0124 / R4 / SCECO / Dead) to the
#   R1 sheet.

$wb = $excel.ActiveWorkbook

# --- R1 -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3946:29:17"
$ws1.Range("G3").Value = "86:01:55"
$ws1.Range("G4").Value = "109:01:55"

# New outage row appended at the bottom of R1
$ws1.Range("B6").Value = "R4"
$ws1.Range("D6").Value = "JED0124"
$ws1.Range("I6").Value = "SCECO"
$ws1.Range("J6").Value = "Dead"
$ws1.Range("L6").Value = "Latis"

# --- R2 -------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12127:52:58"
$ws2.Range("G3").Value = "3257:36:27"
$ws2.Range("G4").Value = "495:48:01"

# --- R4 -------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2973:42:47"
$ws4.Range("G3").Value = "200:55:02"
$ws4.Range("G4").Value = "89:07:27"
$ws4.Range("G5").Value = "86:45:00"

# --- R5 -------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "447:41:46"

# --- R6 -------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "88:14:04"
